$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7853
$ws.Range("F3").Value = 8007
$ws.Range("F6").Value = 7014
$ws.Range("F7").Value = 3462
$ws.Range("F9").Value = 3777
$ws.Range("F12").Value = 62
$ws.Range("F14").Value = 109
$ws.Range("F15").Value = 493
$ws.Range("F17").Value = 98
$ws.Range("F20").Value = 4
$ws.Range("F22").Value = 4001
$ws.Range("F26").Value = 517
$ws.Range("F27").Value = 1590
$ws.Range("F28").Value = 94
$ws.Range("F30").Value = 2907
$ws.Range("F31").Value = 2093
$ws.Range("F32").Value = 51
$ws.Range("F34").Value = 78
$ws.Range("F35").Value = 29
$ws.Range("F36").Value = 20
$ws.Range("F38").Value = 4032
$ws.Range("F39").Value = 405
$ws.Range("F40").Value = 306
$ws.Range("F41").Value = 50
$ws.Range("F42").Value = 937
$ws.Range("F43").Value = 709
$ws.Range("F44").Value = 118
$ws.Range("F45").Value = 1551
$ws.Range("F47").Value = 592
$ws.Range("F48").Value = 686

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 245
$ws.Range("F6").Value = 59
$ws.Range("F7").Value = 106
$ws.Range("F14").Value = 20
$ws.Range("F15").Value = 518

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 245
$ws.Range("F4").Value = 7853
$ws.Range("F5").Value = 8007
$ws.Range("F8").Value = 7015
$ws.Range("F9").Value = 3462
$ws.Range("F11").Value = 3777
$ws.Range("F14").Value = 62
$ws.Range("F15").Value = 109
$ws.Range("F18").Value = 98
$ws.Range("F20").Value = 59
$ws.Range("F21").Value = 106
$ws.Range("F22").Value = 4
$ws.Range("F24").Value = 4001
$ws.Range("F28").Value = 517
$ws.Range("F29").Value = 1590
$ws.Range("F30").Value = 94
$ws.Range("F32").Value = 2907
$ws.Range("F33").Value = 2093
$ws.Range("F34").Value = 51
$ws.Range("F38").Value = 4032
$ws.Range("F40").Value = 405
$ws.Range("F41").Value = 306
$ws.Range("F42").Value = 20
$ws.Range("F43").Value = 50
$ws.Range("F44").Value = 709
$ws.Range("F45").Value = 1551
$ws.Range("F48").Value = 686

Write-Host "Applied F-column (想去人数) updates across 展览, 演出, and 全部类型 sheets."
